# Insert a new "Calizon" distance column before the existing "Poblacion"
# column (old column E), shifting Poblacion and its values to column F,
# then populate the new column with the allocation/report distance data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at E; existing E (Poblacion) shifts to F.
$ws.Columns("E:E").Insert()

# Header for the newly inserted column.
$ws.Range("E1").Value = "Calizon"

# Data values for the new column.
$ws.Range("E2").Value = 1.972318
$ws.Range("E3").Value = 1.87996
$ws.Range("E4").Value = 2.206026
$ws.Range("E5").Value = 4.654394
$ws.Range("E6").Value = 4.654394
